# Daily auto-refresh of the quantum-computing ticker table.
# Date advances to 2025-12-06 and the four tickers are re-ranked / re-scored
# with the day's new metrics (row order becomes IBM, D-Wave, IonQ, Rigetti).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$watch  = [char]0x26D4 + " 관망하십시오."
$neutral = [char]0x26AA + " 중립 구간"

$rows = @(
    @{ Row = 2; Name = "International Business Machines"; Ticker = "IBM";  D = 310.76; E = 54.6; F = 0.71;  G = 60; H = 60; I = 66; J = 63; K = 59.9 },
    @{ Row = 3; Name = "D-Wave Quantum Inc.";              Ticker = "QBTS"; D = 27.85;  E = 62.1; F = 22.85; G = 60; H = 66; I = 66; J = 83; K = 59.9 },
    @{ Row = 4; Name = "IonQ, Inc.";                       Ticker = "IONQ"; D = 53.08;  E = 59.7; F = 7.67;  G = 60; H = 56; I = 60; J = 70; K = 57.5 },
    @{ Row = 5; Name = "Rigetti Computing, Inc.";          Ticker = "RGTI"; D = 28.7;   E = 58.2; F = 12.26; G = 50; H = 60; I = 66; J = 83; K = 56.9 }
)

# Column A holds the date as plain text ("2025-12-05" style), not a real
# date value, so force the text number format before writing the string --
# otherwise a YYYY-MM-DD literal gets auto-parsed into a date serial.
$dateRange = $ws.Range("A2:A5")
$dateRange.NumberFormat = "@"

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = "2025-12-06"
    $ws.Cells.Item($row, 2).Value = $r.Name
    $ws.Cells.Item($row, 3).Value = $r.Ticker
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
    $ws.Cells.Item($row, 8).Value = $r.H
    $ws.Cells.Item($row, 9).Value = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = "Pattern"
    $ws.Cells.Item($row, 13).Value = $watch
    $ws.Cells.Item($row, 14).Value = 51.54219175917372
    $ws.Cells.Item($row, 15).Value = $neutral
}
